$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.404.90'
$ws.Range('E2').Value = '  +0.95%  '
$ws.Range('D3').Value = '1.793.78'
$ws.Range('E3').Value = '  +0.93%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '227.54'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.17%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.555'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +2.11%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '33.03'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +4.69%  '
$ws.Range('E9').Value = '  +1.74%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0692'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +1.40%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0948'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.16%  '
$ws.Range('D12').Value = '2.054.58'
$ws.Range('E12').Value = '  +1.01%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.17'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +3.02%  '
$ws.Range('D14').Value = '1.776.59'
$ws.Range('E14').Value = '  -0.24%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.638'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +3.27%  '
$ws.Range('D16').Value = '34.425.44'
$ws.Range('E16').Value = '  +1.01%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.30'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +3.05%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '68.83'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +2.11%  '
$ws.Range('D19').Value = '0.0₃0803'
$ws.Range('E19').Value = '  +2.38%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '245.90'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.78%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.32'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +2.95%  '
$ws.Range('E22').Value = '  -0.01%  '
$ws.Range('E23').Value = '  +2.48%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '169.98'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +5.86%  '
$ws.Range('E25').Value = '  +1.35%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.38'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +3.95%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.61'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +2.78%  '
$ws.Range('E28').Value = '  +2.04%  '
$ws.Range('E29').Value = '  -0.16%  '
$ws.Range('E30').Value = '  +9.37%  '
$ws.Range('E31').Value = '  +2.26%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.24'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.78%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.81'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +2.92%  '
$ws.Range('E34').Value = '  +2.82%  '
$ws.Range('D35').Value = '1.423.14'
$ws.Range('E35').Value = '  -1.51%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.687'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +5.21%  '
$ws.Range('E37').Value = '  +6.27%  '
$ws.Range('E38').Value = '  +3.39%  '
$ws.Range('E39').Value = '  +0.90%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '84.72'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +6.02%  '
$ws.Range('E41').Value = '  +3.99%  '
$ws.Range('E42').Value = '  +2.25%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.40'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.96%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.03'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +3.74%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0526'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +2.13%  '
$ws.Range('E46').Value = '  +2.89%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '6.15'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +1.64%  '
$ws.Range('D48').Value = '1.955.89'
$ws.Range('E48').Value = '  +1.01%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '105.34'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +1.82%  '
$ws.Range('E50').Value = '  -0.08%  '
$ws.Range('D51').Value = '0.0₆0129'
$ws.Range('E51').Value = '  -0.76%  '
